$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Rename "DUMPLING EATERY" -> "DUMPLING HOUSE EP" (row 18, Customer Name)
$ws.Range("A18").Value = "DUMPLING HOUSE EP"

# Fill in the previously-blank "Last Invoice Date" for WAKAN TIPI CENTER (row 9),
# matching the date formatting used by the other date cells in column D.
$ws.Range("D7").Copy()
$ws.Range("D9").PasteSpecial(-4122)
$ws.Range("D9").Value = 45950

# Insert two new customer rows above the last row (old row 20 / HOLY FAMILY
# MARONITE CHURCH), which pushes that row down to row 22.
$ws.Rows("20:21").Insert()
$ws.Rows("20:21").RowHeight = $ws.Rows("19").RowHeight

# Carry the (blank) formatting of column F down into the newly inserted rows.
$ws.Range("F19").Copy()
$ws.Range("F20:F21").PasteSpecial(-4122)

# New row 20: ELITE LIQUOR INC
$ws.Range("A20").Value = "ELITE LIQUOR INC"
$ws.Range("B20").Value = "Larsen, Rick J"
$ws.Range("C20").Value = "023"
$ws.Range("E20").Value = "0008344"

# New row 21: Executive Aviation
$ws.Range("A21").Value = "Executive Aviation"
$ws.Range("B21").Value = "Dack, Suzanne"
$ws.Range("C21").Value = "015"
$ws.Range("E21").Value = "0008345"
